$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1408
$ws1.Range("F6").Value = 5

# Update "全部类型" sheet (sheetId 4), which mirrors the data above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1408
$ws4.Range("F6").Value = 5
